$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1) ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E3").Value = 572.50464375
$schedule.Range("F3").Value = 18.93203187003968

$schedule.Range("E4").Value = 330.5716245
$schedule.Range("F4").Value = 21.86320267857143

# --- Sheet "Detailed" (sheet2) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B11").Value = 79.95005
$detailed.Range("B12").Value = 79.95

$detailed.Range("C14").Value = "historical"

$detailed.Range("B16").Value = 56.98
$detailed.Range("B18").Value = 56.97999
$detailed.Range("B19").Value = 48.86725
$detailed.Range("B20").Value = 36.06

$detailed.Range("B23").Value = 36.06
$detailed.Range("B24").Value = 36.07

$detailed.Range("B27").Value = 36.0601

$detailed.Range("B32").Value = 28.9852
$detailed.Range("B33").Value = 19.62131
$detailed.Range("B34").Value = 19.23499

$detailed.Range("B38").Value = -3.17461
$detailed.Range("B39").Value = -2.70941
$detailed.Range("B40").Value = 0.01121

$detailed.Range("B42").Value = 29.85322
$detailed.Range("B43").Value = 22.01959
$detailed.Range("B44").Value = 0.04217
$detailed.Range("B45").Value = 59.01628
$detailed.Range("B46").Value = 57.04922
$detailed.Range("B47").Value = 57.04367
$detailed.Range("B48").Value = 57.04367
